$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.422.82'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').Value = '1.850.40'
$ws.Range('E3').Value = '  +1.17%  '
$ws.Range('D4').Value = "'1.000"
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'233.08"
$ws.Range('E5').Value = '  +1.51%  '
$ws.Range('D6').Value = "'1.000"
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').Value = "'0.4764"
$ws.Range('E7').Value = '  +3.06%  '
$ws.Range('D8').Value = "'0.2743"
$ws.Range('E8').Value = '  +1.69%  '
$ws.Range('E9').Value = '  +1.49%  '
$ws.Range('D10').Value = "'17.59"
$ws.Range('E10').Value = '  +9.76%  '
$ws.Range('D11').Value = '1.853.49'
$ws.Range('E11').Value = '  +1.42%  '
$ws.Range('D12').Value = "'0.07466"
$ws.Range('E12').Value = '  +1.41%  '
$ws.Range('D13').Value = "'4.945"
$ws.Range('E13').Value = '  +1.12%  '
$ws.Range('D14').Value = "'84.62"
$ws.Range('E14').Value = '  +2.10%  '
$ws.Range('D15').Value = "'0.6238"
$ws.Range('E15').Value = '  +0.95%  '
$ws.Range('D16').Value = '30.388.39'
$ws.Range('E16').Value = '  +1.29%  '
$ws.Range('D17').Value = "'245.48"
$ws.Range('E17').Value = '  +8.15%  '
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('E19').Value = '  +2.80%  '
$ws.Range('D20').Value = "'0.000007322"
$ws.Range('E20').Value = '  +0.89%  '
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('D22').Value = "'4.903"
$ws.Range('E22').Value = '  +1.69%  '
$ws.Range('D23').Value = "'5.898"
$ws.Range('E23').Value = '  +1.62%  '
$ws.Range('D24').Value = "'164.63"
$ws.Range('E24').Value = '  -0.39%  '
$ws.Range('D25').Value = "'9.077"
$ws.Range('D26').Value = "'17.99"
$ws.Range('E26').Value = '  +1.39%  '
$ws.Range('D27').Value = "'1.868"
$ws.Range('E27').Value = '  +1.79%  '
$ws.Range('E28').Value = '  +1.81%  '
$ws.Range('D29').Value = "'1.351"
$ws.Range('E29').Value = '  -1.12%  '
$ws.Range('D30').Value = "'4.036"
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('D31').Value = "'3.813"
$ws.Range('E31').Value = '  +1.93%  '
$ws.Range('D32').Value = "'0.04826"
$ws.Range('E32').Value = '  +1.01%  '
$ws.Range('D33').Value = "'1.128"
$ws.Range('E33').Value = '  +0.62%  '
$ws.Range('D34').Value = "'0.6944"
$ws.Range('E34').Value = '  -0.34%  '
$ws.Range('D35').Value = "'2.704"
$ws.Range('E35').Value = '  +0.64%  '
$ws.Range('D36').Value = "'0.01895"
$ws.Range('E36').Value = '  +5.28%  '
$ws.Range('D37').Value = "'2.681"
$ws.Range('E37').Value = '  +3.13%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = "'0.8755"
$ws.Range('E38').Value = '  -1.41%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = "'1.993"
$ws.Range('E39').Value = '  +4.55%  '
$ws.Range('D40').Value = "'106.68"
$ws.Range('E40').Value = '  +4.06%  '
$ws.Range('D41').Value = "'1.000"
$ws.Range('E41').Value = '  +0.09%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = "'0.4049"
$ws.Range('E42').Value = '  +1.92%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = "'5.502"
$ws.Range('E43').Value = '  +0.67%  '
$ws.Range('D44').Value = "'7.158"
$ws.Range('E44').Value = '  +3.85%  '
$ws.Range('D45').Value = "'62.91"
$ws.Range('E45').Value = '  +6.48%  '
$ws.Range('D46').Value = "'0.1194"
$ws.Range('E46').Value = '  +0.81%  '
$ws.Range('D47').Value = "'33.67"
$ws.Range('E47').Value = '  +3.78%  '
$ws.Range('D48').Value = "'8.514"
$ws.Range('E48').Value = '  +1.65%  '
$ws.Range('D49').Value = "'0.05505"
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('D50').Value = "'1.347"
$ws.Range('E50').Value = '  +0.21%  '
$ws.Range('D51').Value = "'0.3675"
$ws.Range('E51').Value = '  +1.71%  '
